$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 110 (shifts existing rows 110-161 down to 113-164)
$ws.Rows("110:112").Insert()

# Row 110
$ws.Cells.Item(110, 1).Value = 11
$ws.Cells.Item(110, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(110, 3).Value = 'Bíobío'
$ws.Cells.Item(110, 4).Value = 44596
$ws.Cells.Item(110, 5).Value = 8
$ws.Cells.Item(110, 6).Value = 'Fruta'
$ws.Cells.Item(110, 7).Value = 100103
$ws.Cells.Item(110, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(110, 9).Value = 100103004
$ws.Cells.Item(110, 10).Value = 'Durazno'
$ws.Cells.Item(110, 11).Value = 'Carson'
$ws.Cells.Item(110, 12).Value = 'Primera'
$ws.Cells.Item(110, 13).Value = 220
$ws.Cells.Item(110, 14).Value = 10000
$ws.Cells.Item(110, 15).Value = 11000
$ws.Cells.Item(110, 16).Value = 10545
$ws.Cells.Item(110, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(110, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(110, 19).Value = 659
$ws.Cells.Item(110, 20).Value = 16

# Row 111
$ws.Cells.Item(111, 1).Value = 11
$ws.Cells.Item(111, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(111, 3).Value = 'Bíobío'
$ws.Cells.Item(111, 4).Value = 44596
$ws.Cells.Item(111, 5).Value = 8
$ws.Cells.Item(111, 6).Value = 'Fruta'
$ws.Cells.Item(111, 7).Value = 100103
$ws.Cells.Item(111, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(111, 9).Value = 100103004
$ws.Cells.Item(111, 10).Value = 'Durazno'
$ws.Cells.Item(111, 11).Value = 'Elegant Lady'
$ws.Cells.Item(111, 12).Value = 'Primera'
$ws.Cells.Item(111, 13).Value = 220
$ws.Cells.Item(111, 14).Value = 9000
$ws.Cells.Item(111, 15).Value = 10000
$ws.Cells.Item(111, 16).Value = 9545
$ws.Cells.Item(111, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(111, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(111, 19).Value = 597
$ws.Cells.Item(111, 20).Value = 16

# Row 112
$ws.Cells.Item(112, 1).Value = 11
$ws.Cells.Item(112, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(112, 3).Value = 'Bíobío'
$ws.Cells.Item(112, 4).Value = 44596
$ws.Cells.Item(112, 5).Value = 8
$ws.Cells.Item(112, 6).Value = 'Fruta'
$ws.Cells.Item(112, 7).Value = 100103
$ws.Cells.Item(112, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(112, 9).Value = 100103004
$ws.Cells.Item(112, 10).Value = 'Durazno'
$ws.Cells.Item(112, 11).Value = 'Polar King'
$ws.Cells.Item(112, 12).Value = 'Primera'
$ws.Cells.Item(112, 13).Value = 220
$ws.Cells.Item(112, 14).Value = 9000
$ws.Cells.Item(112, 15).Value = 10000
$ws.Cells.Item(112, 16).Value = 9455
$ws.Cells.Item(112, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(112, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(112, 19).Value = 591
$ws.Cells.Item(112, 20).Value = 16
